$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.495.02'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.922.49'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.57%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.01'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4698'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.37%  '
$ws.Range('E8').Value = '  -3.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06734'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '106.29'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '18.28'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07754'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.906.24'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.89%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.301'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6585'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '290.16'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.491.77'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.99%  '
$ws.Range('E18').Value = '  +0.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007584'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.79%  '
$ws.Range('E20').Value = '  -2.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.148.90'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.257'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.50%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.200'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.378'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '169.23'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.27'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.126'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1067'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.368'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.172'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.57%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.979'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.99%  '
$ws.Range('E33').Value = '  -1.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7397'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.81%  '
$ws.Range('E35').Value = '  -1.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02090'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.65%  '
$ws.Range('E37').Value = '  +0.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.679'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.060'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '110.19'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.70%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8717'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.857'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.36%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4249'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.63%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '67.24'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '49.66'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +16.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.194'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.234'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.05'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1215'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.2460'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +9.97%  '
